# Update constant values on the "constants" worksheet, reflecting the
# latest pickle file loaded from saved_uncertainty_analyses directory.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value  = 23.70558129644467
$ws.Range("B3").Value  = 0.7293812134916696
$ws.Range("B5").Value  = 1900.817683522394
$ws.Range("B6").Value  = 177447.9812439665
$ws.Range("B13").Value = 0.7278559562870588
$ws.Range("B14").Value = 0.6015346188309691
$ws.Range("B15").Value = 2.747031447972637
$ws.Range("B16").Value = 1.317098396709526
$ws.Range("B18").Value = 1933
